# SSU_kupovina_slicice.docx edit
# 1) Heading "2.1 Kratak opis" paragraph: explicit zero left indent
#    <w:ind w:hanging="0"/> -> <w:ind w:left="0" w:hanging="0"/>
$d = $word.ActiveDocument

$foundHeading = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Style.NameLocal -eq "Heading 2" -and $para.Range.Text -like "*Kratak opis*") {
        $para.Format.LeftIndent = 0
        $foundHeading = $true
        break
    }
}
if (-not $foundHeading) {
    throw "Could not locate the '2.1 Kratak opis' Heading 2 paragraph"
}

# 2) Footer: the two anchored shapes (arrow connector + page-number bracket)
#    were nudged by the author (slightly larger extents) and Word re-minted
#    the content-control IDs for the page-number field on both the
#    DrawingML and VML-fallback branches.
$footerXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><w:pPr><w:pStyle w:val="Footer"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor behindDoc="1" distT="0" distB="0" distL="0" distR="0" simplePos="0" locked="0" layoutInCell="1" allowOverlap="1" relativeHeight="6" wp14:anchorId="2380B3D6"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:align>center</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:align>center</wp:align></wp:positionV><wp:extent cx="5520055" cy="2540"/><wp:effectExtent l="9525" t="9525" r="6350" b="9525"/><wp:wrapNone/><wp:docPr id="1" name="Straight Arrow Connector 1"/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5519520" cy="1800"/></a:xfrm><a:custGeom><a:avLst/><a:gdLst/><a:ahLst/><a:rect l="l" t="t" r="r" b="b"/><a:pathLst><a:path w="21600" h="21600"><a:moveTo><a:pt x="0" y="0"/></a:moveTo><a:lnTo><a:pt x="21600" y="21600"/></a:lnTo></a:path></a:pathLst></a:custGeom><a:noFill/><a:ln w="12600"><a:solidFill><a:srgbClr val="808080"/></a:solidFill><a:round/></a:ln></wps:spPr><wps:style><a:lnRef idx="0"></a:lnRef><a:fillRef idx="0"/><a:effectRef idx="0"></a:effectRef><a:fontRef idx="minor"/></wps:style><wps:bodyPr/></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict></w:pict></mc:Fallback></mc:AlternateContent><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor behindDoc="1" distT="0" distB="0" distL="0" distR="0" simplePos="0" locked="0" layoutInCell="1" allowOverlap="1" relativeHeight="11" wp14:anchorId="1A394E62"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:align>center</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:align>center</wp:align></wp:positionV><wp:extent cx="564515" cy="240665"/><wp:effectExtent l="19050" t="19050" r="19685" b="18415"/><wp:wrapNone/><wp:docPr id="2" name="Double Bracket 2"/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="563760" cy="240120"/></a:xfrm><a:prstGeom prst="bracketPair"><a:avLst><a:gd name="adj" fmla="val 16667"/></a:avLst></a:prstGeom><a:solidFill><a:srgbClr val="ffffff"/></a:solidFill><a:ln w="28440"><a:solidFill><a:srgbClr val="808080"/></a:solidFill><a:round/></a:ln></wps:spPr><wps:style><a:lnRef idx="0"></a:lnRef><a:fillRef idx="0"/><a:effectRef idx="0"></a:effectRef><a:fontRef idx="minor"/></wps:style><wps:txbx><w:txbxContent><w:sdt><w:sdtPr><w:docPartObj><w:docPartGallery w:val="Page Numbers (Bottom of Page)"/><w:docPartUnique w:val="true"/></w:docPartObj><w:id w:val="1905599996"/></w:sdtPr><w:sdtContent><w:p><w:pPr><w:pStyle w:val="FrameContents"/><w:spacing w:before="0" w:after="160"/><w:jc w:val="center"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:fldChar w:fldCharType="begin"></w:fldChar></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:instrText> PAGE </w:instrText></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:sdtContent></w:sdt></w:txbxContent></wps:txbx><wps:bodyPr tIns="0" bIns="0"><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>10000</wp14:pctWidth></wp14:sizeRelH></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shapetype id="shapetype_185" coordsize="21600,21600" o:spt="185" adj="3600" path="m0@0qy@6@7l@1,qx@8@6l21600@2qy@9@10l@0,21600qx@7@9xnsem@0,21600qx@7@9l0@0qy@6@7m@1,qx@8@6l21600@2qy@9@10nfe"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="val #0"/><v:f eqn="sum width 0 @0"/><v:f eqn="sum height 0 @0"/><v:f eqn="prod @0 2929 10000"/><v:f eqn="sum width 0 @3"/><v:f eqn="sum height 0 @3"/><v:f eqn="sum @0 0 0"/><v:f eqn="sum 0 @0 @0"/><v:f eqn="sum @0 @1 0"/><v:f eqn="sum 0 21600 @0"/><v:f eqn="sum @0 @2 0"/></v:formulas><v:path gradientshapeok="t" o:connecttype="rect" textboxrect="@3,@3,@4,@5"/><v:handles><v:h position="0,@0"/></v:handles></v:shapetype><v:shape id="shape_0" ID="Double Bracket 2" fillcolor="white" stroked="t" style="position:absolute;margin-left:211.8pt;margin-top:-3.05pt;width:44.35pt;height:18.85pt;mso-position-horizontal:center;mso-position-horizontal-relative:margin;mso-position-vertical:center" wp14:anchorId="1A394E62" type="shapetype_185"><w10:wrap type="square"/><v:fill o:detectmouseclick="t" type="solid" color2="black"/><v:stroke color="gray" weight="28440" joinstyle="round" endcap="flat"/><v:textbox><w:txbxContent><w:sdt><w:sdtPr><w:docPartObj><w:docPartGallery w:val="Page Numbers (Bottom of Page)"/><w:docPartUnique w:val="true"/></w:docPartObj><w:id w:val="1584151582"/></w:sdtPr><w:sdtContent><w:p><w:pPr><w:pStyle w:val="FrameContents"/><w:spacing w:before="0" w:after="160"/><w:jc w:val="center"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:fldChar w:fldCharType="begin"></w:fldChar></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:instrText> PAGE </w:instrText></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:sdtContent></w:sdt></w:txbxContent></v:textbox></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p>
'@

$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.InsertXML($footerXml)
